# docs : ERD 자료 수정
# Fix a handful of incorrect "Physical" (column C) values in the
# RentalReservation ERD worksheet - mostly copy/paste leftovers from
# other tables' PK/FK column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = "HCP_RENTABLE_PRODUCT_RESERVATION_INFO_UID"
$ws.Range("C27").Value = "HCP_RENTABLE_PRODUCT_STOCK_CATEGORY_UID"
$ws.Range("C29").Value = "HCP_RENTABLE_PRODUCT_STOCK_IMAGE_UID"
$ws.Range("C76").Value = "HCP_RENTABLE_PRODUCT_RESERVATION_INFO_UID"
$ws.Range("C97").Value = "HCP_RENTABLE_PRODUCT_STOCK_CATEGORY_UID"
$ws.Range("C105").Value = "HCP_RENTABLE_PRODUCT_STOCK_IMAGE_UID"
$ws.Range("C113").Value = "HCP_PAYMENT_UID"
$ws.Range("C127").Value = "HCP_PAYMENT_UID"
$ws.Range("C137").Value = "HCP_PAYMENT_UID"
$ws.Range("C147").Value = "HCP_PAYMENT_UID"
$ws.Range("C155").Value = "HCP_PAYMENT_UID"
$ws.Range("C162").Value = "HCP_RENTABLE_PRODUCT_RESERVATION_PAYMENT_INFO_UID"
$ws.Range("C165").Value = "HCP_RENTABLE_PRODUCT_RESERVATION_INFO_UID"
$ws.Range("C166").Value = "HCP_PAYMENT_UID"
